# ADDITIONAL SCRAPING: add "ODI Bowling Extra" sheet (MAIDEN_OVERS /
# PERCENT_WICKETS_OF_ALL per match) and drop the now-unused blank
# attribute cells from "ODI Batting Extra".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "ODI Batting Extra" - remove cells that hold no data (columns
#    B..E, BATTING_POSITION/NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL) so only
#    cells that actually carry a value remain.
# ---------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$lastRow = $battingExtra.Cells.Item(1,1).Worksheet.UsedRange.Rows.Count

for ($r = 2; $r -le 21; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $battingExtra.Cells.Item($r, $c)
        if ($cell.Value2 -eq "") {
            $cell.ClearContents()
        }
    }
}

# ---------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet after the last existing
#    sheet, with sheetId/position following "ODI Batting Extra".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $afterSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row - bold, centered, top-aligned, thin border all round
# (matches the look of every other sheet's header row).
$header = $bowlingExtra.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$bowlingExtra.Cells.Item(1,1).NumberFormat = "@"
$bowlingExtra.Cells.Item(1,1).Value = "MATCH_CODE"
$bowlingExtra.Cells.Item(1,2).NumberFormat = "@"
$bowlingExtra.Cells.Item(1,2).Value = "MAIDEN_OVERS"
$bowlingExtra.Cells.Item(1,3).NumberFormat = "@"
$bowlingExtra.Cells.Item(1,3).Value = "PERCENT_WICKETS_OF_ALL"

# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL - one row per match
$rows = @(
     ("3781","1","20.00%")
    ,("3786","","")
    ,("3791","1","")
    ,("3795","0","10.00%")
    ,("3808","0","30.00%")
    ,("3810","","")
    ,("3811","1","20.00%")
    ,("3841","0","10.00%")
    ,("3874","0","20.00%")
    ,("3875","","")
    ,("3974","0","")
    ,("3976","","")
    ,("3978","","")
    ,("4042","0","10.00%")
    ,("4047","0","")
    ,("4050","","")
    ,("4052","0","10.00%")
    ,("4053","","")
    ,("4524","","")
    ,("4526","1","")
)

$r = 2
foreach ($row in $rows) {
    $bowlingExtra.Cells.Item($r, 1).NumberFormat = "@"
    $bowlingExtra.Cells.Item($r, 1).Value = $row[0]
    $bowlingExtra.Cells.Item($r, 2).NumberFormat = "@"
    $bowlingExtra.Cells.Item($r, 2).Value = $row[1]
    $bowlingExtra.Cells.Item($r, 3).NumberFormat = "@"
    $bowlingExtra.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
